$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Sheet1" to "Template"
$ws.Name = "Template"

# Populate header + data cells (set in an order that matches the
# original shared-strings table order: Analyte, CBD, {{ cbd }}, Result)
$ws.Range("B5").Value = "Analyte"
$ws.Range("B8").Value = "CBD"
$ws.Range("C8").Value = "{{ cbd }}"
$ws.Range("C5").Value = "Result"

# Column widths
$ws.Columns.Item(3).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 19

# Final selection on C8
$ws.Range("C8").Select() | Out-Null
